$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.986.94'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.16%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.310.40'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.57%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.06'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.15%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.90'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.95%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  -2.22%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.310.78'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.50%  '

$ws.Range('E10').Value = '  +0.18%  '

$ws.Range('E11').Value = '  +0.63%  '

$ws.Range('E12').Value = '  -0.52%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.332'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.38%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.30'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.50%  '

$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.985.58'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.21%  '

$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.724.59'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.66%  '

$ws.Range('E17').Value = '  -0.72%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.306.17'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.85%  '

$ws.Range('E19').Value = '  -1.51%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.05'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.40%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '311.02'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.13%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.53'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.67%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.33%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.35'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.93%  '

$ws.Range('E25').Value = '  -3.17%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.00%  '

$ws.Range('E27').Value = '  -3.19%  '

$ws.Range('E28').Value = '  +2.16%  '

$ws.Range('B29').Value = 'SuiNetwork'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.20'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.26%  '

$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.81'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.25%  '

$ws.Range('E31').Value = '  -0.78%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0723'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.53%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.83'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.99%  '

$ws.Range('B34').Value = 'PolygonEcosystemToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.378'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.73%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.34'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -6.50%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.65'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.38%  '

$ws.Range('E38').Value = '  +0.05%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.01'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.44%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '317.44'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.08%  '

$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.51'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.20%  '

$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '37.48'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.47%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '136.22'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.47%  '

$ws.Range('E44').Value = '  -0.63%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0944'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.24%  '

$ws.Range('E46').Value = '  +1.79%  '

$ws.Range('E47').Value = '  +1.57%  '

$ws.Range('E48').Value = '  -0.91%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₆0222'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +21.35%  '

$ws.Range('E50').Value = '  +1.07%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.01'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.06%  '
